$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (40_heterogeneous / deterministic) - only Execution Time changes
$ws.Range("H2").Value = 0.002998828887939453

# Row 3 (40_heterogeneous / randomized)
$ws.Range("C3").Value = 627.63
$ws.Range("D3").Value = 489.41
$ws.Range("E3").Value = 138.22
$ws.Range("F3").Value = 1117.04
$ws.Range("G3").Value = 558.52
$ws.Range("H3").Value = 1.749794721603394

# Row 4 (40_homogeneous / deterministic) - only Execution Time changes
$ws.Range("H4").Value = 0.002999305725097656

# Row 5 (40_homogeneous / randomized)
$ws.Range("C5").Value = 514.4
$ws.Range("D5").Value = 514.27
$ws.Range("E5").Value = 0.14
$ws.Range("F5").Value = 1028.67
$ws.Range("G5").Value = 514.34
$ws.Range("H5").Value = 1.640929460525513

# Row 6 (60_homogeneous / deterministic) - only Execution Time changes
$ws.Range("H6").Value = 0.006999492645263672

# Row 7 (60_homogeneous / randomized)
$ws.Range("C7").Value = 555.95
$ws.Range("D7").Value = 553.0700000000001
$ws.Range("E7").Value = 2.88
$ws.Range("F7").Value = 1664.91
$ws.Range("G7").Value = 554.97
$ws.Range("H7").Value = 3.689326286315918

# Row 8 (80_heterogeneous / deterministic) - only Execution Time changes
$ws.Range("H8").Value = 0.01099991798400879

# Row 9 (80_heterogeneous / randomized)
$ws.Range("C9").Value = 690.84
$ws.Range("D9").Value = 469.68
$ws.Range("E9").Value = 221.16
$ws.Range("F9").Value = 2331.17
$ws.Range("G9").Value = 582.79
$ws.Range("H9").Value = 6.423193216323853

# Row 10 (80_homogeneous / deterministic) - only Execution Time changes
$ws.Range("H10").Value = 0.01000046730041504

# Row 11 (80_homogeneous / randomized)
$ws.Range("C11").Value = 537.59
$ws.Range("D11").Value = 530.78
$ws.Range("E11").Value = 6.81
$ws.Range("F11").Value = 2136.53
$ws.Range("G11").Value = 534.13
$ws.Range("H11").Value = 6.032886981964111
